$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 903373.4
$ws.Range("I86").Value = 13155.272
$ws.Range("J86").Value = 1329129.9
$ws.Range("K86").Value = 13155.272
$ws.Range("L86").Value = 1329129.9
$ws.Range("M86").Value = -12032.272
$ws.Range("N86").Value = -1331375.9
$ws.Range("H89").Value = 903373.4
$ws.Range("I89").Value = 13155.272
$ws.Range("J89").Value = 1329129.9
$ws.Range("K89").Value = 65776.36
$ws.Range("L89").Value = 6645649.5
$ws.Range("M89").Value = -60160.36
$ws.Range("N89").Value = -6656881.5
$ws.Range("H121").Value = 909.6875
$ws.Range("J121").Value = 909.6875
$ws.Range("L121").Value = 2729.0625
$ws.Range("N121").Value = -6223.0625
$ws.Range("H131").Value = 1471.0526
$ws.Range("I131").Value = 421.25
$ws.Range("K131").Value = 1263.75
$ws.Range("M131").Value = 3776.25
$ws.Range("H137").Value = 2000.55
$ws.Range("I137").Value = 2037.4878
$ws.Range("J137").Value = 1920.8422
$ws.Range("K137").Value = 6112.463400000001
$ws.Range("L137").Value = 5762.5266
$ws.Range("M137").Value = -3562.463400000001
$ws.Range("N137").Value = -10862.5266
$ws.Range("H141").Value = 1690.9524
$ws.Range("I141").Value = 1184.7368
$ws.Range("J141").Value = 6500
$ws.Range("K141").Value = 3554.2104
$ws.Range("L141").Value = 19500
$ws.Range("M141").Value = 1625.7896
$ws.Range("N141").Value = -29860

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 30406916
$ws.Range("I61").Value = 46876148
$ws.Range("J61").Value = 2180.7693
$ws.Range("K61").Value = 46876148
$ws.Range("L61").Value = 2180.7693
$ws.Range("M61").Value = -46875936
$ws.Range("N61").Value = -2604.7693
$ws.Range("H136").Value = 30406916
$ws.Range("I136").Value = 46876148
$ws.Range("J136").Value = 2180.7693
$ws.Range("K136").Value = 140628444
$ws.Range("L136").Value = 6542.3079
$ws.Range("M136").Value = -140625894
$ws.Range("N136").Value = -11642.3079

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1059269.4
$ws.Range("I86").Value = 2650.5
$ws.Range("J86").Value = 2327212
$ws.Range("K86").Value = 2650.5
$ws.Range("L86").Value = 2327212
$ws.Range("M86").Value = -1527.5
$ws.Range("N86").Value = -2329458
$ws.Range("H89").Value = 1059269.4
$ws.Range("I89").Value = 2650.5
$ws.Range("J89").Value = 2327212
$ws.Range("K89").Value = 13252.5
$ws.Range("L89").Value = 11636060
$ws.Range("M89").Value = -7636.5
$ws.Range("N89").Value = -11647292
$ws.Range("H112").Value = 78800
$ws.Range("J112").Value = 78800
$ws.Range("L112").Value = 78800
$ws.Range("N112").Value = -81754
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1262.3846
$ws.Range("I16").Value = 1664.2
$ws.Range("J16").Value = 1011.25
$ws.Range("K16").Value = 1664.2
$ws.Range("L16").Value = 1011.25
$ws.Range("M16").Value = -1377.2
$ws.Range("N16").Value = -1585.25
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H94").Value = 4396.8076
$ws.Range("I94").Value = 20114.4
$ws.Range("J94").Value = 654.5238000000001
$ws.Range("K94").Value = 20114.4
$ws.Range("L94").Value = 654.5238000000001
$ws.Range("M94").Value = -19663.4
$ws.Range("N94").Value = -1556.5238
$ws.Range("H113").Value = 1262.3846
$ws.Range("I113").Value = 1664.2
$ws.Range("J113").Value = 1011.25
$ws.Range("K113").Value = 1664.2
$ws.Range("L113").Value = 1011.25
$ws.Range("M113").Value = 505.8
$ws.Range("N113").Value = -5351.25
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 30575.777
$ws.Range("I12").Value = 5.428571
$ws.Range("J12").Value = 37954.83
$ws.Range("K12").Value = 16.285713
$ws.Range("L12").Value = 113864.49
$ws.Range("M12").Value = 156.714287
$ws.Range("N12").Value = -114210.49
$ws.Range("H105").Value = 5875
$ws.Range("J105").Value = 5875
$ws.Range("L105").Value = 17625
$ws.Range("N105").Value = -22867
$ws.Range("H129").Value = 1217.1
$ws.Range("I129").Value = 858.2
$ws.Range("J129").Value = 1576
$ws.Range("K129").Value = 2574.6
$ws.Range("L129").Value = 4728
$ws.Range("M129").Value = 2425.4
$ws.Range("N129").Value = -14728
$ws.Range("H131").Value = 806.4400000000001
$ws.Range("I131").Value = 442.23077
$ws.Range("J131").Value = 860.86206
$ws.Range("K131").Value = 1326.69231
$ws.Range("L131").Value = 2582.58618
$ws.Range("M131").Value = 3713.30769
$ws.Range("N131").Value = -12662.58618

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 100000000
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 100000000
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 4331.364
$ws.Range("I132").Value = 2156.7954
$ws.Range("J132").Value = 13029.637
$ws.Range("K132").Value = 6470.3862
$ws.Range("L132").Value = 39088.911
$ws.Range("M132").Value = -3940.3862
$ws.Range("N132").Value = -44148.911

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1583.9166
$ws.Range("I82").Value = 1200.5
$ws.Range("J82").Value = 2350.75
$ws.Range("K82").Value = 1200.5
$ws.Range("L82").Value = 2350.75
$ws.Range("M82").Value = -839.5
$ws.Range("N82").Value = -3072.75
$ws.Range("H85").Value = 1583.9166
$ws.Range("I85").Value = 1200.5
$ws.Range("J85").Value = 2350.75
$ws.Range("K85").Value = 1200.5
$ws.Range("L85").Value = 2350.75
$ws.Range("M85").Value = 47.5
$ws.Range("N85").Value = -4846.75
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H136").Value = 39307110
$ws.Range("I136").Value = 44735516
$ws.Range("J136").Value = 33335870
$ws.Range("K136").Value = 134206548
$ws.Range("L136").Value = 100007610
$ws.Range("M136").Value = -134203998
$ws.Range("N136").Value = -100012710

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 527.5294
$ws.Range("I100").Value = 443.23077
$ws.Range("K100").Value = 886.46154
$ws.Range("M100").Value = -345.46154
$ws.Range("H126").Value = 2531.35
$ws.Range("I126").Value = 2223
$ws.Range("J126").Value = 3250.8333
$ws.Range("K126").Value = 6669
$ws.Range("L126").Value = 9752.499899999999
$ws.Range("M126").Value = -4199
$ws.Range("N126").Value = -14692.4999
